# Fix spelling error on the "Your input" slide:
# "What is the common practice invite in your system ..." ->
# "What is the common practice in your system ..."
# The corrected paragraph is re-authored as three runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(3)

$run1 = $para.Runs(1)

# Rewrite the run with the corrected (spelling-fixed) text first.
$run1.Text = "What is the common practice in your system for the email containing an invite after inserting the event into the users calendar?"

# Split the corrected text into three runs, matching the target authoring:
#   "What is the " | "common practice in " | "your system for the email containing an invite after inserting the event into the users calendar?"
$run3 = $run1.Characters(32, 97)
$run3.Text = "your system for the email containing an invite after inserting the event into the users calendar?"

$run2 = $run1.Characters(13, 19)
$run2.Text = "common practice in "
